$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = "32.767884,34.966961"
$ws.Range("B8").Value = "33.084756,35.112133"
$ws.Range("B10").Value = "32.163217,34.961133"
$ws.Range("B11").Value = "32.163217,34.961133"
$ws.Range("B18").Value = "31.854920,35.218710"
$ws.Range("B19").Value = "32.019045,34.841227"
$ws.Range("B21").Value = "32.163217,34.961133"
$ws.Range("B24").Value = "32.773251,35.044543"
$ws.Range("B25").Value = "32.606459,35.290914"
$ws.Range("B28").Value = "32.707298,35.173307"
$ws.Range("B29").Value = "31.854920,35.218710"
$ws.Range("B32").Value = "32.197179,34.956413"
$ws.Range("B34").Value = "31.767207,35.224441"
$ws.Range("B35").Value = "32.052244,34.797643"
$ws.Range("B36").Value = "32.047035,34.899314"
$ws.Range("B37").Value = "31.790191,35.198620"
$ws.Range("B38").Value = "31.748498,35.214655"
$ws.Range("B41").Value = "32.234686,34.954455"
$ws.Range("B42").Value = "32.197179,34.956413"
$ws.Range("B48").Value = "31.886326,34.786696"
$ws.Range("B49").Value = "32.020531,34.799529"
$ws.Range("B51").Value = "32.037040,34.776415"
$ws.Range("B56").Value = "32.081982,34.816659"
$ws.Range("B57").Value = "32.028209,34.802593"
$ws.Range("B58").Value = "33.132610,35.690627"
$ws.Range("B61").Value = "31.767207,35.224441"
$ws.Range("B64").Value = "32.064156,34.854185"
$ws.Range("B69").Value = "32.151070,34.847113"
$ws.Range("B70").Value = "32.006200,34.743653"
$ws.Range("B75").Value = "32.037040,34.776415"
$ws.Range("B77").Value = "31.818922,35.194455"
$ws.Range("B83").Value = "32.028209,34.802593"
$ws.Range("B85").Value = "31.942541,34.872538"
$ws.Range("B88").Value = "32.097875,34.896155"
$ws.Range("B96").Value = "31.194371,34.837706"
$ws.Range("B97").Value = "31.748475,35.212194"
$ws.Range("B104").Value = "32.174844,34.814576"
$ws.Range("B109").Value = "31.785936,35.221741"
$ws.Range("B111").Value = "31.785077,34.693905"
$ws.Range("B114").Value = "32.163217,34.961133"
$ws.Range("B116").Value = "31.975998,34.882170"
$ws.Range("B119").Value = "31.669726,34.779153"
$ws.Range("B120").Value = "32.090535,34.885715"
$ws.Range("B145").Value = "32.163217,34.961133"
$ws.Range("B184").Value = "32.783385,35.036378"
$ws.Range("B211").Value = "32.163217,34.961133"
$ws.Range("B224").Value = "32.163217,34.961133"
$ws.Range("B243").Value = "32.927663,35.151740"
$ws.Range("B247").Value = "32.777112,35.040416"
$ws.Range("B255").Value = "32.209639,34.964653"
$ws.Range("B258").Value = "31.748475,35.212194"
$ws.Range("B267").Value = "31.750030,35.218672"
$ws.Range("B268").Value = "32.197179,34.956413"
$ws.Range("B270").Value = "32.099723,34.828750"
$ws.Range("B275").Value = "32.047831,34.870851"
$ws.Range("B276").Value = "32.983671,35.251911"
$ws.Range("B282").Value = "31.818001,34.669804"
$ws.Range("B283").Value = "31.667321,34.601532"
$ws.Range("B289").Value = "32.830002,34.970337"
$ws.Range("B295").Value = "31.748475,35.212194"
$ws.Range("B296").Value = "31.785263,35.186530"
$ws.Range("B297").Value = "31.748475,35.212194"
$ws.Range("B299").Value = "31.748475,35.212194"
$ws.Range("B301").Value = "31.748498,35.214655"
$ws.Range("B307").Value = "32.193962,34.884145"
$ws.Range("B308").Value = "31.781986,35.164617"
$ws.Range("B312").Value = "33.003303,35.091469"
$ws.Range("B336").Value = "31.785087,35.210391"
$ws.Range("B339").Value = "31.850998,35.218030"
$ws.Range("B340").Value = "31.757946,35.215351"
$ws.Range("B342").Value = "31.980488,34.813327"
$ws.Range("B347").Value = "32.051929,34.814324"
$ws.Range("B349").Value = "31.893720,34.803882"
$ws.Range("B354").Value = "32.197179,34.956413"
$ws.Range("B356").Value = "32.153195,34.846595"
$ws.Range("B361").Value = "31.746771,35.212895"
$ws.Range("B362").Value = "32.197179,34.956413"
$ws.Range("B365").Value = "31.783743,34.661290"
$ws.Range("B372").Value = "31.755751,34.983774"
$ws.Range("B373").Value = "31.767207,35.224441"
$ws.Range("B376").Value = "32.100317,34.828653"
$ws.Range("B384").Value = "32.163217,34.961133"
$ws.Range("B387").Value = "32.163217,34.961133"
$ws.Range("B388").Value = "32.166162,34.810351"
$ws.Range("B390").Value = "31.992690,34.909264"
$ws.Range("B392").Value = "31.813664,34.667678"
$ws.Range("B398").Value = "32.054678,34.804617"
$ws.Range("B400").Value = "32.011816,34.796993"
$ws.Range("B401").Value = "31.937727,34.837262"
$ws.Range("B403").Value = "31.663407,34.599960"
$ws.Range("B407").Value = "31.675767,34.597809"
$ws.Range("B410").Value = "31.942541,34.872538"
$ws.Range("B415").Value = "31.855315,35.221246"
$ws.Range("B417").Value = "31.750585,35.215673"
$ws.Range("B419").Value = "31.855315,35.221246"
$ws.Range("B420").Value = "31.792434,35.144275"
$ws.Range("B422").Value = "31.755957,34.989832"
$ws.Range("B424").Value = "32.006200,34.743653"
$ws.Range("B427").Value = "31.226237,34.809557"
$ws.Range("B434").Value = "31.860925,35.220821"
$ws.Range("B435").Value = "31.226237,34.809557"
$ws.Range("B441").Value = "31.222301,34.819116"
$ws.Range("B445").Value = "32.167359,34.823770"
$ws.Range("B452").Value = "31.226727,34.804803"
$ws.Range("B453").Value = "31.828293,34.662887"
$ws.Range("B455").Value = "31.666065,34.584567"
$ws.Range("B457").Value = "32.008568,34.807885"
$ws.Range("B458").Value = "32.095945,34.859010"
$ws.Range("B459").Value = "31.222537,34.880798"
$ws.Range("B460").Value = "31.220502,34.815623"
$ws.Range("B464").Value = "31.237729,34.793597"
$ws.Range("B465").Value = "32.159475,34.836814"
$ws.Range("B474").Value = "32.108556,34.939938"
$ws.Range("B475").Value = "31.663407,34.599960"
$ws.Range("B476").Value = "32.163217,34.961133"
$ws.Range("B483").Value = "31.526474,34.596970"
